$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 426 (pushing existing rows 426-512 down to 428-514)
$ws.Rows.Item(426).Resize(2).Insert()

# Common/shared values for both new rows
$mercadoId = 7
$mercado   = "Terminal Hortofrutícola Agro Chillán"
$region    = "Ñuble"
$codreg    = 16
$catId     = 100112003
$categoria = "Ajo"
$variedad  = "Chino"
$calidad   = "Primera"
$origen    = "China"
$kgUnid    = 10
$clasif    = "Hortaliza"

# New row 426: $/caja 10 kilos entry
$r = 426
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = 45173
$ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $catId
$ws.Cells.Item($r, 7).Value = $categoria
$ws.Cells.Item($r, 8).Value = $variedad
$ws.Cells.Item($r, 9).Value = $calidad
$ws.Cells.Item($r, 10).Value = 100
$ws.Cells.Item($r, 11).Value = 21000
$ws.Cells.Item($r, 12).Value = 21000
$ws.Cells.Item($r, 13).Value = 21000
$ws.Cells.Item($r, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item($r, 15).Value = $origen
$ws.Cells.Item($r, 16).Value = 2100
$ws.Cells.Item($r, 17).Value = $kgUnid
$ws.Cells.Item($r, 18).Value = $clasif

# New row 427: $/malla 10 kilos entry
$r = 427
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = 45173
$ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $catId
$ws.Cells.Item($r, 7).Value = $categoria
$ws.Cells.Item($r, 8).Value = $variedad
$ws.Cells.Item($r, 9).Value = $calidad
$ws.Cells.Item($r, 10).Value = 80
$ws.Cells.Item($r, 11).Value = 23000
$ws.Cells.Item($r, 12).Value = 23000
$ws.Cells.Item($r, 13).Value = 23000
$ws.Cells.Item($r, 14).Value = "`$/malla 10 kilos"
$ws.Cells.Item($r, 15).Value = $origen
$ws.Cells.Item($r, 16).Value = 2300
$ws.Cells.Item($r, 17).Value = $kgUnid
$ws.Cells.Item($r, 18).Value = $clasif

Write-Host "Inserted new rows 426 and 427"
